# v03 - added restrictions table and notes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Restrictions table: re-prefix the binary-string formulas in column AI
# with "0b" so the concatenated bit pattern reads as an explicit binary
# literal (AI2 has its own formula; AI3:AI10 share one formula). ---
$ws.Range("AI2").Formula = '="0b"&AF2&AE2&AD2&AC2&AB2&AA2&Z2&Y2&X2'
$ws.Range("AI3:AI10").Formula = '="0b"&AF3&AE3&AD3&AC3&AB3&AA3&Z3&Y3&X3'

# --- Untick the "O" checkbox on row 4 (linked cell O4), which ripples into
# Z4 (=IF(O4,1,0)) and therefore into the AI4 bit-string result. ---
$ws.Range("O4").Value = $false

# Also push the new state into the checkbox shape itself, where supported.
try {
    $ws.Shapes.Item("Check Box 21").ControlFormat.Value = 0
} catch {
}

# --- Hide the helper columns L:V (used for the restriction-table notes) ---
$ws.Range("L1:V1").EntireColumn.Hidden = $true

# --- Update the sheet's stored selection to match the hidden range ---
[void]$ws.Range("L1:V1048576").Select()
